$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.558.56"
$ws.Range("E2").Value = "  -0.06%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.751.19"
$ws.Range("E3").Value = "  -0.32%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "324.17"
$ws.Range("E5").Value = "  -0.12%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.002"
$ws.Range("E6").Value = "  +0.12%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4488"
$ws.Range("E7").Value = "  +0.11%  "
$ws.Range("E8").Value = "  -1.93%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07453"
$ws.Range("E9").Value = "  -1.03%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "41.32"
$ws.Range("E10").Value = "  -1.92%  "
$ws.Range("E11").Value = "  -2.27%  "
$ws.Range("E12").Value = "  +0.09%  "
$ws.Range("E13").Value = "  -0.38%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.970"
$ws.Range("E14").Value = "  -1.48%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.143"
$ws.Range("E15").Value = "  -1.28%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.753.36"
$ws.Range("E16").Value = "  -0.26%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "93.63"
$ws.Range("E17").Value = "  +0.67%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001055"
$ws.Range("E18").Value = "  -1.02%  "
$ws.Range("E19").Value = "  -0.89%  "
$ws.Range("E20").Value = "  +0.08%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.06"
$ws.Range("E21").Value = "  -0.46%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.730"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "27.605.28"
$ws.Range("E23").Value = "  -0.05%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.17"
$ws.Range("E24").Value = "  -0.92%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.085"
$ws.Range("E25").Value = "  -0.41%  "
$ws.Range("E26").Value = "  +1.52%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.12"
$ws.Range("E27").Value = "  -1.88%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.953.54"
$ws.Range("E28").Value = "  -0.24%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.102"
$ws.Range("E29").Value = "  -1.74%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "125.13"
$ws.Range("E30").Value = "  -0.74%  "
$ws.Range("E31").Value = "  -1.36%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09167"
$ws.Range("E32").Value = "  +0.80%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.653"
$ws.Range("E33").Value = "  +0.60%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.489"
$ws.Range("E34").Value = "  -1.40%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.02282"
$ws.Range("E35").Value = "  -1.11%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "11.70"
$ws.Range("E36").Value = "  -3.97%  "
$ws.Range("E37").Value = "  -0.87%  "
$ws.Range("E38").Value = "  -0.12%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.6267"
$ws.Range("E39").Value = "  -2.40%  "
$ws.Range("E40").Value = "  -0.51%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.177"
$ws.Range("E41").Value = "  -1.63%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.396"
$ws.Range("E42").Value = "  -0.09%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "7.764"
$ws.Range("E43").Value = "  -1.27%  "
$ws.Range("E44").Value = "  -1.67%  "
$ws.Range("E45").Value = "  +0.06%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5852"
$ws.Range("E46").Value = "  -1.19%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "121.76"
$ws.Range("E47").Value = "  -0.35%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.930"
$ws.Range("E48").Value = "  -2.27%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.06880"
$ws.Range("E49").Value = "  +0.01%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.126"
$ws.Range("E50").Value = "  -3.39%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.126"
$ws.Range("E51").Value = "  -0.51%  "
